$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the product placeholder to the new aggregated products string.
$ws.Range("F2").Value = '${record.productNamesString}'

# Update the active selection on the sheet.
$ws.Range("A3").Select()

# Update the workbook window size.
$wb.Windows.Item(1).Width = 16740
$wb.Windows.Item(1).Height = 3660
